$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7 values, following the same pattern as rows 4-6
$ws.Range("B7").Value = (New-Object DateTime(2020, 1, 12))
$ws.Range("B7").NumberFormat = "mm/dd/yyyy"

$ws.Range("D7").Value = "D:\Documents\App development\featherlook\featherlook1.0.py"

$ws.Range("E7").Value = "Major development: script is set up to run from executable on Windows. `nMinor:This version includes a popup that shows the absolute path for the temporary file I make to store the search word so I don't have to use .get() all over the place."
$ws.Range("E7").WrapText = $true

$ws.Range("F7").Value = "1)remove popup for path; 2)add scroll; 3) limit size of lists so that users don't select directory with large search and run takes forever - this is just a toy; 4)describe function of app so its obvious for user what this is. Maybe add popup or labels to guide"

$ws.Range("F8").Select()
